# Applies the "cryptos list" update (Wed Apr 24 04:26:34 UTC 2024) to the
# active worksheet. Updates Price (D) / Volume(1h) (E) figures, and swaps
# three pairs of rows whose coin ordering changed (17/18, 39/40, 45/46).
#
# Price values that look like plain numbers (e.g. "604.97") are prefixed
# with a leading apostrophe so Excel keeps them as text, matching the
# original workbook where every Price cell is an inline/shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2;  Col=4; Val="'66.696.61"},
    @{Row=2;  Col=5; Val="  +0.41%  "},

    @{Row=3;  Col=4; Val="'3.245.66"},
    @{Row=3;  Col=5; Val="  +1.77%  "},

    @{Row=4;  Col=5; Val="  -0.13%  "},

    @{Row=5;  Col=4; Val="'604.97"},
    @{Row=5;  Col=5; Val="  +0.22%  "},

    @{Row=6;  Col=4; Val="'157.40"},
    @{Row=6;  Col=5; Val="  +0.81%  "},

    @{Row=7;  Col=4; Val="'0.999"},
    @{Row=7;  Col=5; Val="  -0.01%  "},

    @{Row=8;  Col=4; Val="'3.244.64"},
    @{Row=8;  Col=5; Val="  +1.75%  "},

    @{Row=9;  Col=4; Val="'0.549"},
    @{Row=9;  Col=5; Val="  -0.11%  "},

    @{Row=10; Col=5; Val="  +2.18%  "},

    @{Row=11; Col=4; Val="'5.80"},
    @{Row=11; Col=5; Val="  -1.00%  "},

    @{Row=12; Col=4; Val="'0.502"},
    @{Row=12; Col=5; Val="  -1.79%  "},

    @{Row=13; Col=5; Val="  +2.96%  "},

    @{Row=14; Col=4; Val="'39.11"},
    @{Row=14; Col=5; Val="  +0.35%  "},

    @{Row=15; Col=4; Val="'3.774.86"},
    @{Row=15; Col=5; Val="  +1.63%  "},

    @{Row=16; Col=4; Val="'66.674.48"},
    @{Row=16; Col=5; Val="  +0.22%  "},

    # Rows 17/18 swap: WrappedEther <-> Polkadot
    @{Row=17; Col=2; Val="Polkadot"},
    @{Row=17; Col=3; Val="https://coinranking.com/coin/25W7FG7om+polkadot-dot"},
    @{Row=17; Col=4; Val="'7.31"},
    @{Row=17; Col=5; Val="  -1.31%  "},

    @{Row=18; Col=2; Val="WrappedEther"},
    @{Row=18; Col=3; Val="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"},
    @{Row=18; Col=4; Val="'3.219.52"},
    @{Row=18; Col=5; Val="  +0.93%  "},

    @{Row=19; Col=4; Val="'0.114"},
    @{Row=19; Col=5; Val="  +1.83%  "},

    @{Row=20; Col=4; Val="'508.30"},
    @{Row=20; Col=5; Val="  -1.02%  "},

    @{Row=21; Col=4; Val="'15.35"},
    @{Row=21; Col=5; Val="  -0.96%  "},

    @{Row=22; Col=4; Val="'0.745"},
    @{Row=22; Col=5; Val="  +1.48%  "},

    @{Row=23; Col=4; Val="'8.07"},
    @{Row=23; Col=5; Val="  -0.91%  "},

    @{Row=24; Col=4; Val="'14.67"},
    @{Row=24; Col=5; Val="  -1.65%  "},

    @{Row=25; Col=4; Val="'86.20"},
    @{Row=25; Col=5; Val="  +1.75%  "},

    @{Row=26; Col=4; Val="'0.176"},
    @{Row=26; Col=5; Val="  +95.65%  "},

    @{Row=27; Col=5; Val="  +0.03%  "},

    @{Row=28; Col=5; Val="  +0.20%  "},

    @{Row=29; Col=4; Val="'9.10"},
    @{Row=29; Col=5; Val="  -0.96%  "},

    @{Row=30; Col=4; Val="'2.35"},
    @{Row=30; Col=5; Val="  -1.55%  "},

    @{Row=31; Col=4; Val="'2.90"},
    @{Row=31; Col=5; Val="  -5.54%  "},

    @{Row=32; Col=4; Val="'6.92"},
    @{Row=32; Col=5; Val="  -1.43%  "},

    @{Row=33; Col=4; Val="'28.28"},
    @{Row=33; Col=5; Val="  +0.58%  "},

    @{Row=34; Col=5; Val="  +0.04%  "},

    @{Row=35; Col=5; Val="  -4.69%  "},

    @{Row=36; Col=4; Val="'6.37"},
    @{Row=36; Col=5; Val="  -2.82%  "},

    @{Row=37; Col=4; Val="0.0₃0807"},
    @{Row=37; Col=5; Val="  +19.47%  "},

    @{Row=38; Col=4; Val="'55.35"},
    @{Row=38; Col=5; Val="  +1.07%  "},

    # Rows 39/40 swap: Bittensor <-> dogwifhat
    @{Row=39; Col=2; Val="dogwifhat"},
    @{Row=39; Col=3; Val="https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"},
    @{Row=39; Col=4; Val="'3.29"},
    @{Row=39; Col=5; Val="  +15.08%  "},

    @{Row=40; Col=2; Val="Bittensor"},
    @{Row=40; Col=3; Val="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"},
    @{Row=40; Col=4; Val="'494.03"},
    @{Row=40; Col=5; Val="  -3.27%  "},

    @{Row=41; Col=4; Val="'0.0424"},
    @{Row=41; Col=5; Val="  -0.28%  "},

    @{Row=42; Col=5; Val="  +2.08%  "},

    @{Row=43; Col=4; Val="'8.75"},
    @{Row=43; Col=5; Val="  -1.25%  "},

    @{Row=44; Col=4; Val="'0.294"},
    @{Row=44; Col=5; Val="  -2.62%  "},

    # Rows 45/46 swap: Maker <-> Fetch.AI
    @{Row=45; Col=2; Val="Fetch.AI"},
    @{Row=45; Col=3; Val="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"},
    @{Row=45; Col=4; Val="'2.47"},
    @{Row=45; Col=5; Val="  +1.41%  "},

    @{Row=46; Col=2; Val="Maker"},
    @{Row=46; Col=3; Val="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"},
    @{Row=46; Col=4; Val="'2.947.60"},
    @{Row=46; Col=5; Val="  +3.24%  "},

    @{Row=47; Col=4; Val="'28.30"},
    @{Row=47; Col=5; Val="  -0.56%  "},

    @{Row=48; Col=4; Val="'2.42"},
    @{Row=48; Col=5; Val="  -0.34%  "},

    @{Row=49; Col=5; Val="  +1.70%  "},

    @{Row=50; Col=5; Val="  -0.03%  "},

    @{Row=51; Col=5; Val="  -1.61%  "}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val
}
